$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$ws.Range("A58").Value = 70
$ws.Range("B58").Value = 146.41
$ws.Range("C58").Value = 12
$ws.Range("D58").Value = 300
$ws.Range("E58").Value = 2000
$ws.Range("F58").Value = 'AWKCPJ DTŚZXL ĄBĆEĘF GHIŁMN ŃOÓQRS UVYŹŻ|'
$ws.Range("G58").Value = -1403.88365943578
$ws.Range("H58").Value = 'DTŚZXL ĄBĆEĘF GHIŁMN ŃOÓQRS UVYŹŻ| AWKCPJ'
$ws.Range("I58").Value = -1403.8837

$ws.Range("A59").Value = 173
$ws.Range("B59").Value = 404.22
$ws.Range("C59").Value = 12
$ws.Range("D59").Value = 300
$ws.Range("E59").Value = 2000
$ws.Range("F59").Value = 'PĄMNAÓ SFXDJŹ BCĆEĘG HIKLŁŃ OQRŚTU VWYZŻ|'
$ws.Range("G59").Value = -1403.88365943578
$ws.Range("H59").Value = 'PĄMNAÓ SFXDJŹ BCĆEĘG HIKLŁŃ OQRŚTU VWYZŻ|'
$ws.Range("I59").Value = -1403.8837

$ws.Range("A60").Value = 37
$ws.Range("B60").Value = 76.48999999999999
$ws.Range("C60").Value = 12
$ws.Range("D60").Value = 300
$ws.Range("E60").Value = 2000
$ws.Range("F60").Value = '|ÓLDXŹ JRTCPZ AĄBĆEĘ FGHIKŁ MNŃOQS ŚUVWYŻ'
$ws.Range("G60").Value = -1403.88365943578
$ws.Range("H60").Value = 'XÓLD|Ź JRTCPZ AĄBĆEĘ FGHIKŁ MNŃOQS ŚUVWYŻ'
$ws.Range("I60").Value = -1426.2392

$ws.Range("A61").Value = 48
$ws.Range("B61").Value = 99.76000000000001
$ws.Range("C61").Value = 12
$ws.Range("D61").Value = 300
$ws.Range("E61").Value = 2000
$ws.Range("F61").Value = 'ÓHCABM ZUŁĘŻL ĄĆDEFG IJKNŃO PQRSŚT VWXYŹ|'
$ws.Range("G61").Value = -1403.88365943578
$ws.Range("H61").Value = 'ÓHCABM ZUŁĘŻL ĄĆDEFG IJKNŃO PQRSŚT VWXYŹ|'
$ws.Range("I61").Value = -1403.8837

$ws.Range("A62").Value = 19
$ws.Range("B62").Value = 50.22
$ws.Range("C62").Value = 12
$ws.Range("D62").Value = 300
$ws.Range("E62").Value = 2000
$ws.Range("F62").Value = 'VFLDŻ| ZYSHTB AĄCĆEĘ GIJKŁM NŃOÓPQ RŚUWXŹ'
$ws.Range("G62").Value = -1403.88365943578
$ws.Range("H62").Value = 'VFLDŻ| ZYSHTB AĄCĆEĘ GIJKŁM NŃOÓPQ RŚUWXŹ'
$ws.Range("I62").Value = -1403.8837

$ws.Range("A63").Value = 31
$ws.Range("B63").Value = 73.48999999999999
$ws.Range("C63").Value = 12
$ws.Range("D63").Value = 300
$ws.Range("E63").Value = 2000
$ws.Range("F63").Value = 'KWNŻT| GŁAEIŃ ĄBCĆDĘ FHJLMO ÓPQRSŚ UVXYZŹ'
$ws.Range("G63").Value = -1403.88365943578
$ws.Range("H63").Value = 'KWNŻT| GŁAEIŃ ĄBCĆDĘ FHJLMO ÓPQRSŚ UVXYZŹ'
$ws.Range("I63").Value = -1403.8837

$ws.Range("A64").Value = 187
$ws.Range("B64").Value = 444.19
$ws.Range("C64").Value = 12
$ws.Range("D64").Value = 300
$ws.Range("E64").Value = 2000
$ws.Range("F64").Value = 'GU|ORK ŹYMLZĄ ABCĆDE ĘFHIJŁ NŃÓPQS ŚTVWXŻ'
$ws.Range("G64").Value = -1403.88365943578
$ws.Range("H64").Value = 'NŃÓPQS ŚTVWXŻ GU|ORK ŹYMLZĄ ABCĆDE ĘFHIJŁ'
$ws.Range("I64").Value = -1403.8837

$ws.Range("A65").Value = 263
$ws.Range("B65").Value = 'Attempt failed!'

$ws.Range("A66").Value = 263
$ws.Range("B66").Value = 600.0599999999999
$ws.Range("C66").Value = 12
$ws.Range("D66").Value = 300
$ws.Range("E66").Value = 2000
$ws.Range("F66").Value = 'CYPSŚŃ JVENIŹ AĄBĆDĘ FGHKLŁ MOÓQRT UWXZŻ|'
$ws.Range("G66").Value = -1403.88365943578
$ws.Range("H66").Value = 'TCSFWH GÓDKĆP UOYMAŻ VINXJE ŃRLŁBŹ ZĄŚĘQ|'
$ws.Range("I66").Value = -1964.7275

$ws.Range("A67").Value = 29
$ws.Range("B67").Value = 68.83
$ws.Range("C67").Value = 12
$ws.Range("D67").Value = 300
$ws.Range("E67").Value = 2000
$ws.Range("F67").Value = 'ŹWNTHĘ ŃĆJMPŚ AĄBCDE FGIKLŁ OÓQRSU VXYZŻ|'
$ws.Range("G67").Value = -1403.88365943578
$ws.Range("H67").Value = 'ŹWNTHĘ ŃĆJMPŚ AĄBCDE FGIKLŁ OÓQRSU VXYZŻ|'
$ws.Range("I67").Value = -1403.8837

$ws.Range("A68").Value = 13
$ws.Range("B68").Value = 30.94
$ws.Range("C68").Value = 12
$ws.Range("D68").Value = 300
$ws.Range("E68").Value = 2000
$ws.Range("F68").Value = 'DIBJVP ŚĘŻKHS AĄCĆEF GLŁMNŃ OÓQRTU WXYZŹ|'
$ws.Range("G68").Value = -1403.88365943578
$ws.Range("H68").Value = 'DIBJVP ŚĘŻKHS AĄCĆEF GLŁMNŃ OÓQRTU WXYZŹ|'
$ws.Range("I68").Value = -1403.8837

$ws.Range("A69").Value = 18
$ws.Range("B69").Value = 43.4
$ws.Range("C69").Value = 12
$ws.Range("D69").Value = 300
$ws.Range("E69").Value = 2000
$ws.Range("F69").Value = 'LVBŁ|S TUŃOŚĘ AĄCĆDE FGHIJK MNÓPQR WXYZŹŻ'
$ws.Range("G69").Value = -1403.88365943578
$ws.Range("H69").Value = 'LVBŁ|S TUŃOŚĘ AĄCĆDE FGHIJK MNÓPQR WXYZŹŻ'
$ws.Range("I69").Value = -1403.8837

$ws.Range("A70").Value = 53
$ws.Range("B70").Value = 115.36
$ws.Range("C70").Value = 12
$ws.Range("D70").Value = 300
$ws.Range("E70").Value = 2000
$ws.Range("F70").Value = 'TSÓDXM Q|JŃIH AĄBCĆE ĘFGKLŁ NOPRŚU VWYZŹŻ'
$ws.Range("G70").Value = -1403.88365943578
$ws.Range("H70").Value = 'TSÓDXM Q|JŃIH AĄBCĆE ĘFGKLŁ NOPRŚU VWYZŹŻ'
$ws.Range("I70").Value = -1403.8837

$ws.Range("A71").Value = 11
$ws.Range("B71").Value = 27.83
$ws.Range("C71").Value = 12
$ws.Range("D71").Value = 300
$ws.Range("E71").Value = 2000
$ws.Range("F71").Value = 'UŻĘŃNC EKHŁFY AĄBĆDG IJLMOÓ PQRSŚT VWXZŹ|'
$ws.Range("G71").Value = -1403.88365943578
$ws.Range("H71").Value = 'UŻHŃNC EKĘŁFY AĄBĆDG IJLMOÓ PQRSŚT VWXZŹ|'
$ws.Range("I71").Value = -1466.0466

$ws.Range("A72").Value = 45
$ws.Range("B72").Value = 91.67
$ws.Range("C72").Value = 12
$ws.Range("D72").Value = 300
$ws.Range("E72").Value = 2000
$ws.Range("F72").Value = 'ŁŻIĘWC MQYSDX AĄBĆEF GHJKLN ŃOÓPRŚ TUVZŹ|'
$ws.Range("G72").Value = -1403.88365943578
$ws.Range("H72").Value = 'ŁŻIĘWC MQYSDX AĄBĆEF GHJKLN ŃOÓPRŚ TUVZŹ|'
$ws.Range("I72").Value = -1403.8837

$ws.Range("A73").Value = 53
$ws.Range("B73").Value = 118.46
$ws.Range("C73").Value = 12
$ws.Range("D73").Value = 300
$ws.Range("E73").Value = 2000
$ws.Range("F73").Value = 'WJDĄCP FVSRKŁ ABĆEĘG HILMNŃ OÓQŚTU XYZŹŻ|'
$ws.Range("G73").Value = -1403.88365943578
$ws.Range("H73").Value = 'WJDĄCP FVSRKŁ ABĆEĘG HILMNŃ OÓQŚTU XYZŹŻ|'
$ws.Range("I73").Value = -1403.8837

$ws.Range("A74").Value = 259
$ws.Range("B74").Value = 600.27
$ws.Range("C74").Value = 12
$ws.Range("D74").Value = 300
$ws.Range("E74").Value = 2000
$ws.Range("F74").Value = 'WIGNTŃ CXJQBO AĄĆDEĘ FHKLŁM ÓPRSŚU VYZŹŻ|'
$ws.Range("G74").Value = -1403.88365943578
$ws.Range("H74").Value = 'WIGNTŃ CXJQBO AĄĆDEĘ FHKLŁM ÓPRSŚU VYZŹŻ|'
$ws.Range("I74").Value = -1403.8837

$ws.Range("A75").Value = 44
$ws.Range("B75").Value = 95.36
$ws.Range("C75").Value = 12
$ws.Range("D75").Value = 300
$ws.Range("E75").Value = 2000
$ws.Range("F75").Value = 'ÓŃTSAB ŻĘELŚJ ĄCĆDFG HIKŁMN OPQRUV WXYZŹ|'
$ws.Range("G75").Value = -1403.88365943578
$ws.Range("H75").Value = 'ÓŃTSAB ŻĘELŚJ ĄCĆDFG HIKŁMN OPQRUV WXYZŹ|'
$ws.Range("I75").Value = -1403.8837

$ws.Range("A76").Value = 41
$ws.Range("B76").Value = 93.18000000000001
$ws.Range("C76").Value = 12
$ws.Range("D76").Value = 300
$ws.Range("E76").Value = 2000
$ws.Range("F76").Value = 'EPBVDŹ NIGŻŁC AĄĆĘFH JKLMŃO ÓQRSŚT UWXYZ|'
$ws.Range("G76").Value = -1403.88365943578
$ws.Range("H76").Value = 'EPBVDŹ NIGŻŁC AĄĆĘFH JKLMŃO ÓQRSŚT UWXYZ|'
$ws.Range("I76").Value = -1403.8837

$ws.Range("A77").Value = 70
$ws.Range("B77").Value = 161.05
$ws.Range("C77").Value = 12
$ws.Range("D77").Value = 300
$ws.Range("E77").Value = 2000
$ws.Range("F77").Value = 'EÓPMTR |ŁJŻWŃ AĄBCĆD ĘFGHIK LNOQSŚ UVXYZŹ'
$ws.Range("G77").Value = -1403.88365943578
$ws.Range("H77").Value = 'EÓPMTR |ŁJŻWŃ AĄBCĆD ĘFGHIK LNOQSŚ UVXYZŹ'
$ws.Range("I77").Value = -1403.8837

$ws.Range("A78").Value = 51
$ws.Range("B78").Value = 116.09
$ws.Range("C78").Value = 12
$ws.Range("D78").Value = 300
$ws.Range("E78").Value = 2000
$ws.Range("F78").Value = 'ZUKMOÓ WPJŻBV AĄCĆDE ĘFGHIL ŁNŃQRS ŚTXYŹ|'
$ws.Range("G78").Value = -1403.88365943578
$ws.Range("H78").Value = 'ZUKMOÓ WPJŻBV AĄCĆDE ĘFGHIL ŁNŃQRS ŚTXYŹ|'
$ws.Range("I78").Value = -1403.8837

$ws.Range("A79").Value = 47
$ws.Range("B79").Value = 110.56
$ws.Range("C79").Value = 12
$ws.Range("D79").Value = 300
$ws.Range("E79").Value = 2000
$ws.Range("F79").Value = 'NMJKGY HŚUÓLS AĄBCĆD EĘFIŁŃ OPQRTV WXZŹŻ|'
$ws.Range("G79").Value = -1403.88365943578
$ws.Range("H79").Value = 'NMJKGY HŚUÓLS AĄBCĆD EĘFIŁŃ OPQRTV WXZŹŻ|'
$ws.Range("I79").Value = -1403.8837

$ws.Range("A80").Value = 46
$ws.Range("B80").Value = 111.47
$ws.Range("C80").Value = 12
$ws.Range("D80").Value = 300
$ws.Range("E80").Value = 2000
$ws.Range("F80").Value = 'EATLHF MVDNQ| ĄBCĆĘG IJKŁŃO ÓPRSŚU WXYZŹŻ'
$ws.Range("G80").Value = -1403.88365943578
$ws.Range("H80").Value = 'EATLHF MVDNQ| ĄBCĆĘG IJKŁŃO ÓPRSŚU WXYZŹŻ'
$ws.Range("I80").Value = -1403.8837

$ws.Range("A81").Value = 258
$ws.Range("B81").Value = 614.6799999999999
$ws.Range("C81").Value = 12
$ws.Range("D81").Value = 300
$ws.Range("E81").Value = 2000
$ws.Range("F81").Value = 'ŃIŁMEO RŚŻYĘG AĄBCĆD FHJKLN ÓPQSTU VWXZŹ|'
$ws.Range("G81").Value = -1403.88365943578
$ws.Range("H81").Value = 'FHJKLN ÓPQSTU VWXZŹ| ŃIŁMEO RŚŻYĘG AĄBCĆD'
$ws.Range("I81").Value = -1403.8837

$ws.Range("A82").Value = 63
$ws.Range("B82").Value = 135.78
$ws.Range("C82").Value = 12
$ws.Range("D82").Value = 300
$ws.Range("E82").Value = 2000
$ws.Range("F82").Value = 'PEVBŃS WLGQHŁ AĄCĆDĘ FIJKMN OÓRŚTU XYZŹŻ|'
$ws.Range("G82").Value = -1403.88365943578
$ws.Range("H82").Value = 'PEVBŃS WLGQHŁ AĄCĆDĘ FIJKMN OÓRŚTU XYZŹŻ|'
$ws.Range("I82").Value = -1403.8837

$ws.Range("A83").Value = 58
$ws.Range("B83").Value = 128.24
$ws.Range("C83").Value = 12
$ws.Range("D83").Value = 300
$ws.Range("E83").Value = 2000
$ws.Range("F83").Value = 'ÓĘĆBUI ŚDLFVX AĄCEGH JKŁMNŃ OPQRST WYZŹŻ|'
$ws.Range("G83").Value = -1403.88365943578
$ws.Range("H83").Value = 'ŚDLXVF AĄCEGH JKŁMNŃ OPQRST WYZŹŻ| ÓĘĆBUI'
$ws.Range("I83").Value = -1447.0448

$ws.Range("A84").Value = 215
$ws.Range("B84").Value = 481.13
$ws.Range("C84").Value = 12
$ws.Range("D84").Value = 300
$ws.Range("E84").Value = 2000
$ws.Range("F84").Value = 'CWEVŚB OŹTNQD AĄĆĘFG HIJKLŁ MŃÓPRS UXYZŻ|'
$ws.Range("G84").Value = -1403.88365943578
$ws.Range("H84").Value = 'CWEVŚB OŹTNQD AĄĆĘFG HIJKLŁ MŃÓPRS UXYZŻ|'
$ws.Range("I84").Value = -1403.8837

$ws.Range("A85").Value = 125
$ws.Range("B85").Value = 262.61
$ws.Range("C85").Value = 12
$ws.Range("D85").Value = 300
$ws.Range("E85").Value = 2000
$ws.Range("F85").Value = 'ŃZKTLF SJ|ŁĘY AĄBCĆD EGHIMN OÓPQRŚ UVWXŹŻ'
$ws.Range("G85").Value = -1403.88365943578
$ws.Range("H85").Value = 'ŃZKTLF SJ|ŁĘY AĄBCĆD EGHIMN OÓPQRŚ UVWXŹŻ'
$ws.Range("I85").Value = -1403.8837

$ws.Range("A86").Value = 31
$ws.Range("B86").Value = 67.7
$ws.Range("C86").Value = 12
$ws.Range("D86").Value = 300
$ws.Range("E86").Value = 2000
$ws.Range("F86").Value = 'JDGLIM ZNRFPE AĄBCĆĘ HKŁŃOÓ QSŚTUV WXYŹŻ|'
$ws.Range("G86").Value = -1403.88365943578
$ws.Range("H86").Value = 'JDGLIM ZNRFPE AĄBCĆĘ HKŁŃOÓ QSŚTUV WXYŹŻ|'
$ws.Range("I86").Value = -1403.8837

$ws.Range("A87").Value = 113
$ws.Range("B87").Value = 247.5
$ws.Range("C87").Value = 12
$ws.Range("D87").Value = 300
$ws.Range("E87").Value = 2000
$ws.Range("F87").Value = 'GWNĘŻJ EMFBYO AĄCĆDH IKLŁŃÓ PQRSŚT UVXZŹ|'
$ws.Range("G87").Value = -1403.88365943578
$ws.Range("H87").Value = 'IKLŁŃÓ PQRSŚT UVXZŹ| GWNĘŻJ EMFBYO AĄCĆDH'
$ws.Range("I87").Value = -1403.8837
